# Applies betexplorer.com re-scrape update to 2023/wales_cymru-premier_2023-2024.xlsx (Sheet1)
# - Rows 75-87: reconciles match order/details for matches played on the same date
#   (F:V columns refreshed in place; A-E index/date columns are untouched)
# - Row 89: appends newly scraped match "Aberystwyth vs Bala" (16-17/11/2023 odds)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75
$ws.Cells.Item(75, 6).Value = "Penybont"
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = "Aberystwyth"
$ws.Cells.Item(75, 9).Value = 2
$ws.Cells.Item(75, 10).Value = 1.26
$ws.Cells.Item(75, 11).Value = "26/10/2023 09:13"
$ws.Cells.Item(75, 12).Value = 1.46
$ws.Cells.Item(75, 13).Value = "27/10/2023 20:36"
$ws.Cells.Item(75, 14).Value = 5.37
$ws.Cells.Item(75, 15).Value = "26/10/2023 09:13"
$ws.Cells.Item(75, 16).Value = 4.44
$ws.Cells.Item(75, 17).Value = "27/10/2023 20:36"
$ws.Cells.Item(75, 18).Value = 8.199999999999999
$ws.Cells.Item(75, 19).Value = "26/10/2023 09:13"
$ws.Cells.Item(75, 20).Value = 6.86
$ws.Cells.Item(75, 21).Value = "27/10/2023 20:36"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/penybont-aberystwyth/EsDH7Voj/"

# Row 76
$ws.Cells.Item(76, 6).Value = "Connahs Q."
$ws.Cells.Item(76, 7).Value = 6
$ws.Cells.Item(76, 8).Value = "Caernarfon"
$ws.Cells.Item(76, 9).Value = 1
$ws.Cells.Item(76, 10).Value = 1.36
$ws.Cells.Item(76, 11).Value = "26/10/2023 09:13"
$ws.Cells.Item(76, 12).Value = 1.38
$ws.Cells.Item(76, 13).Value = "27/10/2023 20:36"
$ws.Cells.Item(76, 14).Value = 4.78
$ws.Cells.Item(76, 15).Value = "26/10/2023 09:13"
$ws.Cells.Item(76, 16).Value = 5.04
$ws.Cells.Item(76, 17).Value = "27/10/2023 20:41"
$ws.Cells.Item(76, 18).Value = 6.43
$ws.Cells.Item(76, 19).Value = "26/10/2023 09:13"
$ws.Cells.Item(76, 20).Value = 7.64
$ws.Cells.Item(76, 21).Value = "27/10/2023 20:41"
$ws.Cells.Item(76, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/connahs-q-caernarfon/KKDD8BWq/"

# Row 77
$ws.Cells.Item(77, 6).Value = "Pontypridd"
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = "Cardiff Metropolitan"
$ws.Cells.Item(77, 9).Value = 3
$ws.Cells.Item(77, 10).Value = 2.65
$ws.Cells.Item(77, 11).Value = "26/10/2023 14:42"
$ws.Cells.Item(77, 12).Value = 3.2
$ws.Cells.Item(77, 13).Value = "28/10/2023 15:21"
$ws.Cells.Item(77, 14).Value = 3.02
$ws.Cells.Item(77, 15).Value = "26/10/2023 14:42"
$ws.Cells.Item(77, 16).Value = 3
$ws.Cells.Item(77, 17).Value = "28/10/2023 15:21"
$ws.Cells.Item(77, 18).Value = 2.61
$ws.Cells.Item(77, 19).Value = "26/10/2023 14:42"
$ws.Cells.Item(77, 20).Value = 2.43
$ws.Cells.Item(77, 21).Value = "28/10/2023 15:21"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/pontypridd-united-cardiff-metropolitan-university/QwG4TC89/"

# Row 78
$ws.Cells.Item(78, 6).Value = "Colwyn Bay"
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = "Barry"
$ws.Cells.Item(78, 9).Value = 1
$ws.Cells.Item(78, 10).Value = 2.28
$ws.Cells.Item(78, 11).Value = "26/10/2023 14:42"
$ws.Cells.Item(78, 12).Value = 2.36
$ws.Cells.Item(78, 13).Value = "28/10/2023 15:22"
$ws.Cells.Item(78, 14).Value = 3.35
$ws.Cells.Item(78, 15).Value = "26/10/2023 14:42"
$ws.Cells.Item(78, 16).Value = 3.71
$ws.Cells.Item(78, 17).Value = "28/10/2023 15:22"
$ws.Cells.Item(78, 18).Value = 2.75
$ws.Cells.Item(78, 19).Value = "26/10/2023 14:42"
$ws.Cells.Item(78, 20).Value = 2.75
$ws.Cells.Item(78, 21).Value = "28/10/2023 15:22"
$ws.Cells.Item(78, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/colwyn-bay-barry-town/hQHdVYwc/"

# Row 79
$ws.Cells.Item(79, 6).Value = "Newtown"
$ws.Cells.Item(79, 7).Value = 1
$ws.Cells.Item(79, 8).Value = "Haverfordwest"
$ws.Cells.Item(79, 9).Value = 1
$ws.Cells.Item(79, 10).Value = 1.81
$ws.Cells.Item(79, 11).Value = "26/10/2023 14:42"
$ws.Cells.Item(79, 12).Value = 1.66
$ws.Cells.Item(79, 13).Value = "28/10/2023 15:28"
$ws.Cells.Item(79, 14).Value = 3.62
$ws.Cells.Item(79, 15).Value = "26/10/2023 14:42"
$ws.Cells.Item(79, 16).Value = 3.91
$ws.Cells.Item(79, 17).Value = "28/10/2023 15:28"
$ws.Cells.Item(79, 18).Value = 3.78
$ws.Cells.Item(79, 19).Value = "26/10/2023 14:42"
$ws.Cells.Item(79, 20).Value = 5
$ws.Cells.Item(79, 21).Value = "28/10/2023 15:28"
$ws.Cells.Item(79, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/newtown-haverfordwest/WGG0Uhg3/"

# Row 80
$ws.Cells.Item(80, 6).Value = "TNS"
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = "Bala"
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 1.12
$ws.Cells.Item(80, 11).Value = "27/10/2023 06:42"
$ws.Cells.Item(80, 12).Value = 1.11
$ws.Cells.Item(80, 13).Value = "28/10/2023 18:13"
$ws.Cells.Item(80, 14).Value = 8.220000000000001
$ws.Cells.Item(80, 15).Value = "27/10/2023 06:42"
$ws.Cells.Item(80, 16).Value = 9.59
$ws.Cells.Item(80, 17).Value = "28/10/2023 18:13"
$ws.Cells.Item(80, 18).Value = 12.17
$ws.Cells.Item(80, 19).Value = "27/10/2023 06:42"
$ws.Cells.Item(80, 20).Value = 19.55
$ws.Cells.Item(80, 21).Value = "28/10/2023 18:13"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/tns-bala/hYBL6k1d/"

# Row 81
$ws.Cells.Item(81, 6).Value = "TNS"
$ws.Cells.Item(81, 7).Value = 6
$ws.Cells.Item(81, 8).Value = "Colwyn Bay"
$ws.Cells.Item(81, 9).Value = 1
$ws.Cells.Item(81, 10).Value = 1.03
$ws.Cells.Item(81, 11).Value = "30/10/2023 16:42"
$ws.Cells.Item(81, 12).Value = 1.02
$ws.Cells.Item(81, 13).Value = "31/10/2023 19:58"
$ws.Cells.Item(81, 14).Value = 16.84
$ws.Cells.Item(81, 15).Value = "30/10/2023 16:42"
$ws.Cells.Item(81, 16).Value = 27.01
$ws.Cells.Item(81, 17).Value = "31/10/2023 20:29"
$ws.Cells.Item(81, 18).Value = 22.13
$ws.Cells.Item(81, 19).Value = "30/10/2023 16:42"
$ws.Cells.Item(81, 20).Value = 47.31
$ws.Cells.Item(81, 21).Value = "31/10/2023 20:29"
$ws.Cells.Item(81, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/tns-colwyn-bay/Y9MMIefi/"

# Row 82
$ws.Cells.Item(82, 6).Value = "Bala"
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = "Pontypridd"
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 1.6
$ws.Cells.Item(82, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(82, 12).Value = 1.57
$ws.Cells.Item(82, 13).Value = "04/11/2023 15:22"
$ws.Cells.Item(82, 14).Value = 3.75
$ws.Cells.Item(82, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(82, 16).Value = 3.68
$ws.Cells.Item(82, 17).Value = "04/11/2023 15:22"
$ws.Cells.Item(82, 18).Value = 4.72
$ws.Cells.Item(82, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(82, 20).Value = 6.74
$ws.Cells.Item(82, 21).Value = "04/11/2023 15:22"
$ws.Cells.Item(82, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/bala-pontypridd-united/vsK8SWNF/"

# Row 83
$ws.Cells.Item(83, 6).Value = "Barry"
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = "Aberystwyth"
$ws.Cells.Item(83, 9).Value = 1
$ws.Cells.Item(83, 10).Value = 1.74
$ws.Cells.Item(83, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(83, 12).Value = 1.69
$ws.Cells.Item(83, 13).Value = "04/11/2023 14:58"
$ws.Cells.Item(83, 14).Value = 3.68
$ws.Cells.Item(83, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(83, 16).Value = 4.07
$ws.Cells.Item(83, 17).Value = "04/11/2023 14:58"
$ws.Cells.Item(83, 18).Value = 3.91
$ws.Cells.Item(83, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(83, 20).Value = 4.5
$ws.Cells.Item(83, 21).Value = "04/11/2023 14:58"
$ws.Cells.Item(83, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/barry-town-aberystwyth/M59DRjwM/"

# Row 84
$ws.Cells.Item(84, 6).Value = "Caernarfon"
$ws.Cells.Item(84, 7).Value = 2
$ws.Cells.Item(84, 8).Value = "Penybont"
$ws.Cells.Item(84, 9).Value = 4
$ws.Cells.Item(84, 10).Value = 2.62
$ws.Cells.Item(84, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(84, 12).Value = 2.62
$ws.Cells.Item(84, 13).Value = "04/11/2023 15:22"
$ws.Cells.Item(84, 14).Value = 3.34
$ws.Cells.Item(84, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(84, 16).Value = 3.77
$ws.Cells.Item(84, 17).Value = "04/11/2023 15:22"
$ws.Cells.Item(84, 18).Value = 2.38
$ws.Cells.Item(84, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(84, 20).Value = 2.44
$ws.Cells.Item(84, 21).Value = "04/11/2023 15:22"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/caernarfon-penybont/0d8HQAhS/"

# Row 85
$ws.Cells.Item(85, 6).Value = "Cardiff Metropolitan"
$ws.Cells.Item(85, 7).Value = 3
$ws.Cells.Item(85, 8).Value = "Connahs Q."
$ws.Cells.Item(85, 9).Value = 1
$ws.Cells.Item(85, 10).Value = 4.39
$ws.Cells.Item(85, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(85, 12).Value = 5.35
$ws.Cells.Item(85, 13).Value = "04/11/2023 15:23"
$ws.Cells.Item(85, 14).Value = 3.87
$ws.Cells.Item(85, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(85, 16).Value = 4.1
$ws.Cells.Item(85, 17).Value = "04/11/2023 15:23"
$ws.Cells.Item(85, 18).Value = 1.62
$ws.Cells.Item(85, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(85, 20).Value = 1.6
$ws.Cells.Item(85, 21).Value = "04/11/2023 15:23"
$ws.Cells.Item(85, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/cardiff-metropolitan-university-connahs-q/hl0qMU0q/"

# Row 86
$ws.Cells.Item(86, 6).Value = "Haverfordwest"
$ws.Cells.Item(86, 7).Value = 5
$ws.Cells.Item(86, 8).Value = "Colwyn Bay"
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 1.69
$ws.Cells.Item(86, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(86, 12).Value = 1.83
$ws.Cells.Item(86, 13).Value = "04/11/2023 15:22"
$ws.Cells.Item(86, 14).Value = 3.75
$ws.Cells.Item(86, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(86, 16).Value = 3.64
$ws.Cells.Item(86, 17).Value = "04/11/2023 15:24"
$ws.Cells.Item(86, 18).Value = 4.07
$ws.Cells.Item(86, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(86, 20).Value = 4.26
$ws.Cells.Item(86, 21).Value = "04/11/2023 15:24"
$ws.Cells.Item(86, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/haverfordwest-colwyn-bay/WvamLlGk/"

# Row 87
$ws.Cells.Item(87, 6).Value = "Newtown"
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = "TNS"
$ws.Cells.Item(87, 9).Value = 2
$ws.Cells.Item(87, 10).Value = 7.29
$ws.Cells.Item(87, 11).Value = "02/11/2023 08:13"
$ws.Cells.Item(87, 12).Value = 8.68
$ws.Cells.Item(87, 13).Value = "04/11/2023 15:13"
$ws.Cells.Item(87, 14).Value = 5.68
$ws.Cells.Item(87, 15).Value = "02/11/2023 08:13"
$ws.Cells.Item(87, 16).Value = 5.96
$ws.Cells.Item(87, 17).Value = "04/11/2023 15:13"
$ws.Cells.Item(87, 18).Value = 1.27
$ws.Cells.Item(87, 19).Value = "02/11/2023 08:13"
$ws.Cells.Item(87, 20).Value = 1.29
$ws.Cells.Item(87, 21).Value = "04/11/2023 15:10"
$ws.Cells.Item(87, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/newtown-tns/E3lhK8Ve/"

# New row 89 (append new match result)
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = "wales"
$ws.Cells.Item(89, 3).Value = "cymru-premier"
$ws.Cells.Item(89, 4).Value = "2023-2024"
$ws.Cells.Item(89, 5).Value = 45247.875
$ws.Cells.Item(89, 6).Value = "Aberystwyth"
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = "Bala"
$ws.Cells.Item(89, 9).Value = 3
$ws.Cells.Item(89, 10).Value = 4.6
$ws.Cells.Item(89, 11).Value = "16/11/2023 09:13"
$ws.Cells.Item(89, 12).Value = 4.61
$ws.Cells.Item(89, 13).Value = "17/11/2023 20:50"
$ws.Cells.Item(89, 14).Value = 3.75
$ws.Cells.Item(89, 15).Value = "16/11/2023 09:13"
$ws.Cells.Item(89, 16).Value = 3.31
$ws.Cells.Item(89, 17).Value = "17/11/2023 20:50"
$ws.Cells.Item(89, 18).Value = 1.62
$ws.Cells.Item(89, 19).Value = "16/11/2023 09:13"
$ws.Cells.Item(89, 20).Value = 1.85
$ws.Cells.Item(89, 21).Value = "17/11/2023 20:50"
$ws.Cells.Item(89, 22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/aberystwyth-bala/QwhdJSp2/"

# Copy cell formatting (bold/bordered index style, date-time number format) so row 89 matches the rest of the sheet
$ws.Range("A88").Copy()
$ws.Range("A89").PasteSpecial(-4122)
$ws.Range("E88").Copy()
$ws.Range("E89").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
